{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// The document is a one-page essay (title, author, email, blank line, body,\n// \"Summary\" heading, summary paragraph). The commit swaps the whole essay\n// from a \"history\" theme to a \"math & music\" theme: title, author name,\n// e-mail address, and every sentence of body text are replaced, a couple of\n// sentences are merged/split differently, and a trailing blank paragraph is\n// appended at the very end of the document.\n//\n// Because every run inside a given paragraph shares identical formatting\n// (same rFonts / color / sz), sentence-level search & replace is safe here:\n// it reproduces the same visible text/paragraph structure as the reference\n// edit without needing to mirror the exact (arbitrary) run-splitting.\n\nconst replacements = [\n  [\"Echoes of the Past: History's Lessons for the Present\", \"Decoding the Interplay between Math and Music: Harmony Revealed\"],\n  [\"Oliver Walsh\", \"Evelyn Monroe\"],\n  [\"oliver\", \"e\"],\n  [\"walsh@academicmail\", \"monroe@schoolmail\"],\n  [\"Our collective journey through the annals of history holds a mirror to the present, revealing patterns, parallels, and poignant echoes that resonate across time\", \"In the realm of human artistry, music and mathematics have long captivated individuals with their ability to evoke emotion and reveal hidden patterns\"],\n  [\" History is not merely a chronological record of events; it is a tapestry interwoven with lessons, both stark and subtle, that offer invaluable guidance as we navigate the complexities of the present\", \" While they may, on the surface, appear disparate, a closer examination reveals an intriguing interplay that elevates both disciplines\"],\n  [\" From the grand sweep of empires to the intimate dramas of human lives, history whispers secrets that can illuminate our own choices and help us avoid repeating past mistakes\", \" From the rhythmic sequences of beats to the harmonious blend of notes, mathematical principles permeate the very essence of music, shaping its structure and enhancing its aesthetic appeal\"],\n  [\"Historians, like detectives following a trail of clues, meticulously piece together the fragments of bygone eras, uncovering truths that may have been obscured by the mists of time\", \"At the heart of this harmonious alliance lies the concept of ratios\"],\n  [\" They delve into ancient scrolls, decipher cryptic inscriptions, and examine relics that bear witness to the triumphs and travails of civilizations long gone\", \" Similar to the numerical world, musical melodies follow mathematical ratios in their intervals, frequencies, and overtones\"],\n  [\" Through their painstaking efforts, they resurrect lost worlds and breathe life into forgotten characters, enabling us to learn from the wisdom and folly of those who came before us\", \" Each note's unique character stems from its subtle relationship to the other notes it harmonizes with, mirroring mathematical patterns observed in equations and functions. Furthermore, the intrinsic math of musical structure, elucidated in time signatures, scales, and chords, constructs a foundation enabling musicians to craft enchanting compositions\"],\n  [\"History, like a wise elder, imparts sagacity through its tales of past glories and bitter defeats\", \"Venturing deeper into the musical tapestry, one encounters the captivating concept of the Fibonacci sequence\"],\n  [\" It reminds us that the relentless march of time can both heal wounds and reveal vulnerabilities\", \" This renowned sequence, wherein each number equals the sum of the two preceding it, epitomizes nature's inherent mathematical principles\"],\n  [\" From the ruins of fallen empires, we glean insights into the dangers of hubris and the corrosive nature of power\", \" From leaf arrangements to the patterns of seashells, the Fibonacci sequence transcends artistic boundaries, revealing harmony in numerous phenomena\"],\n  [\" The struggles of marginalized groups echo in our contemporary battles for justice and equality. Triumph over adversity becomes a beacon of hope, inspiring us to persevere in the face of challenges\", \" Its manifestation in music becomes evident in the spiraling melodies and harmonious progressions beloved by composers across genres\"],\n  [\"Our journey into the realm of history is a pursuit of wisdom, a search for lessons that can illuminate our paths in the present\", \"Mathematics and music intertwine in a symphony of patterns and principles, unveiling a hidden unity that enriches both domains\"],\n  [\" From the rise and fall of civilizations to the triumphs and struggles of individuals, history's echoes reverberate with meaning and relevance\", \" From note ratios and scales to the structure of musical compositions, mathematical foundations provide the framework on which musical masterpieces are constructed\"],\n  [\" It challenges us to learn from past mistakes, to recognize patterns, and to strive for a future that honors the legacy of those who came before us\", \" As music charms our ears, it also challenges our minds, revealing an elegant collaboration between two worlds often perceived as separate\"],\n  [\" The lessons imparted by history are invaluable, timeless, \", \" In this intricate dance, \"],\n  [\"and eternally relevant, guiding our choices and inspiring us to shape a world worthy of our shared heritage\", \"mathematics plays the conductor, guiding musical expressions into harmonious melodies that touch our hearts and souls\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  // Replace only the first (and expected-only) match, keeping its formatting.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// A new, empty paragraph is appended at the very end of the document body.\ncontext.document.body.insertParagraph(\"\", \"End\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# The document is a one-page essay (title, author, email, blank line, body,\n# \"Summary\" heading, summary paragraph). The commit swaps the whole essay\n# from a \"history\" theme to a \"math & music\" theme: title, author name,\n# e-mail address, and every sentence of body text are replaced, a couple of\n# sentences are merged/split differently, and a trailing blank paragraph is\n# appended at the very end of the document.\n#\n# Because every run inside a given paragraph shares identical formatting\n# (same rFonts / color / sz), sentence-level Find & Replace is safe here: it\n# reproduces the same visible text/paragraph structure as the reference edit\n# without needing to mirror the exact (arbitrary) run-splitting.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n\n# Plain \"Find then set Range.Text\" (rather than Find.Execute's own\n# Replacement parameter) is used deliberately: going through\n# Find.Execute(..., ReplaceWith, ReplaceAll) triggers Word's \"smart quotes\"\n# AutoCorrect and silently turns straight apostrophes (') into curly ones\n# (\\u2019), which would corrupt words like \"note's\"/\"nature's\". Setting\n# Range.Text directly performs a literal replacement and keeps the run's\n# existing formatting (font/color/size) intact.\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Wrap = $wdFindContinue\n    $result = $find.Execute()\n    if (-not $result) {\n        throw \"Text not found: $old\"\n    }\n    $find.Parent.Text = $new\n}\n\nReplace-Text \"Echoes of the Past: History's Lessons for the Present\" \"Decoding the Interplay between Math and Music: Harmony Revealed\"\nReplace-Text \"Oliver Walsh\" \"Evelyn Monroe\"\nReplace-Text \"oliver\" \"e\"\nReplace-Text \"walsh@academicmail\" \"monroe@schoolmail\"\n\nReplace-Text \"Our collective journey through the annals of history holds a mirror to the present, revealing patterns, parallels, and poignant echoes that resonate across time\" \"In the realm of human artistry, music and mathematics have long captivated individuals with their ability to evoke emotion and reveal hidden patterns\"\nReplace-Text \" History is not merely a chronological record of events; it is a tapestry interwoven with lessons, both stark and subtle, that offer invaluable guidance as we navigate the complexities of the present\" \" While they may, on the surface, appear disparate, a closer examination reveals an intriguing interplay that elevates both disciplines\"\nReplace-Text \" From the grand sweep of empires to the intimate dramas of human lives, history whispers secrets that can illuminate our own choices and help us avoid repeating past mistakes\" \" From the rhythmic sequences of beats to the harmonious blend of notes, mathematical principles permeate the very essence of music, shaping its structure and enhancing its aesthetic appeal\"\n\nReplace-Text \"Historians, like detectives following a trail of clues, meticulously piece together the fragments of bygone eras, uncovering truths that may have been obscured by the mists of time\" \"At the heart of this harmonious alliance lies the concept of ratios\"\nReplace-Text \" They delve into ancient scrolls, decipher cryptic inscriptions, and examine relics that bear witness to the triumphs and travails of civilizations long gone\" \" Similar to the numerical world, musical melodies follow mathematical ratios in their intervals, frequencies, and overtones\"\nReplace-Text \" Through their painstaking efforts, they resurrect lost worlds and breathe life into forgotten characters, enabling us to learn from the wisdom and folly of those who came before us\" \" Each note's unique character stems from its subtle relationship to the other notes it harmonizes with, mirroring mathematical patterns observed in equations and functions. Furthermore, the intrinsic math of musical structure, elucidated in time signatures, scales, and chords, constructs a foundation enabling musicians to craft enchanting compositions\"\n\nReplace-Text \"History, like a wise elder, imparts sagacity through its tales of past glories and bitter defeats\" \"Venturing deeper into the musical tapestry, one encounters the captivating concept of the Fibonacci sequence\"\nReplace-Text \" It reminds us that the relentless march of time can both heal wounds and reveal vulnerabilities\" \" This renowned sequence, wherein each number equals the sum of the two preceding it, epitomizes nature's inherent mathematical principles\"\nReplace-Text \" From the ruins of fallen empires, we glean insights into the dangers of hubris and the corrosive nature of power\" \" From leaf arrangements to the patterns of seashells, the Fibonacci sequence transcends artistic boundaries, revealing harmony in numerous phenomena\"\nReplace-Text \" The struggles of marginalized groups echo in our contemporary battles for justice and equality. Triumph over adversity becomes a beacon of hope, inspiring us to persevere in the face of challenges\" \" Its manifestation in music becomes evident in the spiraling melodies and harmonious progressions beloved by composers across genres\"\n\nReplace-Text \"Our journey into the realm of history is a pursuit of wisdom, a search for lessons that can illuminate our paths in the present\" \"Mathematics and music intertwine in a symphony of patterns and principles, unveiling a hidden unity that enriches both domains\"\nReplace-Text \" From the rise and fall of civilizations to the triumphs and struggles of individuals, history's echoes reverberate with meaning and relevance\" \" From note ratios and scales to the structure of musical compositions, mathematical foundations provide the framework on which musical masterpieces are constructed\"\nReplace-Text \" It challenges us to learn from past mistakes, to recognize patterns, and to strive for a future that honors the legacy of those who came before us\" \" As music charms our ears, it also challenges our minds, revealing an elegant collaboration between two worlds often perceived as separate\"\nReplace-Text \" The lessons imparted by history are invaluable, timeless, \" \" In this intricate dance, \"\nReplace-Text \"and eternally relevant, guiding our choices and inspiring us to shape a world worthy of our shared heritage\" \"mathematics plays the conductor, guiding musical expressions into harmonious melodies that touch our hearts and souls\"\n\n# A new, empty paragraph is appended at the very end of the document.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n\n$d.Save()\n"}
